# Generate Report for Handback
# Refreshes the timestamp values recorded on the handback-status report:
#   - Overview sheet: "Latest HO Xliff Generate Date" for the first file (row 2)
#   - zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for row 2
#   - de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for row 2
#
# Note: Overview!G2 and de-de!H2 originally shared the same timestamp
# (2016-08-26 21:03:25) and both move to the same new timestamp
# (2016-08-26 21:04:19), so both cells are updated explicitly.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-26 21:04:19"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-26 21:04:14"
$wsZhCn.Range("K2").Value = "2016-08-26 21:04:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-26 21:04:19"
$wsDeDe.Range("K2").Value = "2016-08-26 21:04:39"
